# edit.ps1
# Applies the MADCS worksheet content update described by the commit:
# refactors stop/time/verification sections, adds new fields for
# blocking/non-blocking averias, order quantities and finalization flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Bulk cell value updates (existing + newly introduced cells).
#    Many of these keep identical visible text but the underlying
#    shared-string catalogue changed; a handful contain genuinely new
#    text per the commit message (averia bloqueante/no bloqueante,
#    parar trabajo, nuevas cantidades, finalizar orden, etc.).
# ---------------------------------------------------------------------
$values = @{
  "C2" = 'Pestañas'
  "C3" = 'Verificacio inicia'
  "D3" = 'Teimpo'
  "E3" = ' calidad'
  "F3" = 'Comsumo'
  "G3" = 'Salidas'
  "C6" = 'Verificacio inicia'
  "C8" = 'Producto'
  "D8" = 'Lote'
  "E8" = 'Cantidad'
  "F8" = 'Verficado'
  "B9" = 'Filtradio por olp'
  "C9" = 'a'
  "D9" = 5
  "E9" = 0
  "F9" = 'Marcado/Desmarcado'
  "C10" = 'a'
  "D10" = 4
  "E10" = 6
  "C11" = 'a'
  "D11" = 8
  "E11" = 3
  "C14" = 'Si todo no tieien valor marcado no pueden hacer nada ni tiempo, consumos y salidas'
  "C18" = ' calidad'
  "C19" = 'Parametros maquinaria'
  "C20" = 'Enlave con excell'
  "C21" = 'Foto'
  "C22" = 'Nº ruta'
  "C23" = 'Recoger muestra'
  "D23" = 'Si/No--sistema No permita cerra /registrar salidas si tiene que coger muestra y no esta marcada la verficacion'
  "C27" = 'Consumos'
  "C29" = 'CR (Consumir por restps)'
  "D29" = 'Quien sirve opl'
  "C30" = 'Si'
  "D30" = 'Almacen'
  "C31" = 'No'
  "D31" = 'Fabrica'
  "C38" = 'Boton de cponsumir Quien sirce Fabrica'
  "E38" = 'Solo aploca a estosd'
  "C43" = 'pq el resto (quien srirve opl  toma valor almacen o fabrica cuando se emplee opcion anterior)hay que hacer uno a uno para hacer los restos de cada lote'
  "E46" = 'Producto'
  "F46" = 'Lote'
  "G46" = 'Cantidad resto'
  "E47" = 'a'
  "F47" = 5
  "G47" = 0
  "C48" = 'Filtrado por opl'
  "E48" = 'a'
  "F48" = 4
  "G48" = 6
  "E49" = 'a'
  "F49" = 8
  "G49" = 3
  "C54" = 'Modificar columnas'
  "C55" = 'Q origina'
  "C56" = 'Q servida '
  "C57" = 'Q restos'
  "C62" = 'Tiempos'
  "E62" = 'Tiempos'
  "D63" = 'USUARIO'
  "E63" = 'ESTADO'
  "F63" = 'tipo averia'
  "D64" = 12
  "E64" = 'Limpieza'
  "D65" = 15
  "E65" = 'Preparación'
  "E66" = 'Avería no bloqueante'
  "F66" = 'según el código de paro'
  "C68" = 'prepra'
  "D68" = 'ejecu'
  "E68" = 'limpieza'
  "F68" = 'averia bloqueante'
  "G68" = 'averia no bloqueante'
  "D71" = 'parar trabajo'
  "C74" = 'Cada vez que se inicia una operación se para la que tenga activa anteriormente'
  "C75" = 'Mismo usuario no puede tener 2 activades simultaneamente de estas'
  "C76" = 'Hay dos tipos de avería, bloqueante  y se para todo y otras no bloqueante pero que permiten seguir funcionando'
  "C77" = 'Usuario que marca la opcion de averia si es bloqueante cierra el resto de tiempos al resto de usuarios y solo apareceria el'
  "C78" = 'Fijar hora ytrabajo palnta de forma que si se esta fuera del horario de trabajo(ojo este cambia según temporada) se cierran todos los tiempos'
  "C80" = 'boton para averias bloqueantes'
  "C81" = 'boton para averias no bloqueantes'
  "D82" = 'incidencia'
  "F82" = 'Esta parte es inmformativa'
  "C83" = 'tipo incidencia'
  "C89" = 'Salidas'
  "E89" = 'Tiempos'
  "C90" = 'Producto'
  "D90" = 'Cantidad a fabricar'
  "E90" = 'Cantidad fabricada'
  "F90" = 'Cantidad pendiente'
  "C92" = 'Dar salida'
  "D92" = 'Finalizar orden'
  "C94" = 'Hay que tener algún sistema para indicar que la orden ya puede ser registrada, que está terminada'
  "C95" = 'Al dar a finalizar orden verificar que los consumos están hechos, que los tiempos están imputados, que las salidas están hechas y parar todos los tiempos pendientes y marcar la orden como finalizada'
}

foreach ($addr in $values.Keys) {
  $ws.Range($addr).Value2 = $values[$addr]
}

# E71 must start with "->" ; Excel stores this with a quote-prefix style
# (text forced, not a formula) so we assign it through .Formula with a
# leading apostrophe, which reproduces the quotePrefix="1" cell style.
$ws.Range("E71").Formula = "'->Para mi trabajo activo e imputa los tiempos a la orden"

# ---------------------------------------------------------------------
# 2) Remove the old C82 cell entirely (its content moved to new C83).
# ---------------------------------------------------------------------
$ws.Range("C82").Clear()

# ---------------------------------------------------------------------
# 3) Formatting for newly introduced cells, matching the established
#    palette already used elsewhere on the sheet:
#      - red fill   (FFFF0000) -> "prepra/ejecu/limpieza/averia.." row
#      - orange fill(FFFFC000) -> "tipo incidencia" row
#      - blue fill + centered  -> section header bars
# ---------------------------------------------------------------------
$red    = 255        # BGR for FFFF0000
$orange = 49407       # BGR for FFFFC000
$blue   = 15773696    # BGR for FF00B0F0

$ws.Range("G68").Interior.Color = $red

$ws.Range("C83").Interior.Color = $orange

$headerRange = $ws.Range("C89:F89")
$headerRange.Interior.Color = $blue
$headerRange.HorizontalAlignment = -4108

$ws.Range("C92:D92").Interior.Color = $red

# ---------------------------------------------------------------------
# 4) New merged header cell for the "Salidas / Tiempos" block.
# ---------------------------------------------------------------------
$ws.Range("C89:F89").Merge()

# ---------------------------------------------------------------------
# 5) Update the view state: scroll down and select C96 like the saved
#    file did.
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("C96").Select()
